$wb = $excel.ActiveWorkbook

# --- Update the two renamed labels ("OS-Drive" -> "my OS-drive", "Data-Drive" -> "my data-drive") ---
# These already exist (with the old text) in row 1 of the fieldnames, URL and comments sheets.
$labelSheets = @(1, 2, 4)
foreach ($idx in $labelSheets) {
    $ws = $wb.Worksheets.Item($idx)
    $ws.Range("A1").Value = "my OS-drive"
    $ws.Range("B1").Value = "my data-drive"
}

# --- Reset the saved selection on those sheets to A1 ---
foreach ($idx in $labelSheets) {
    $wb.Worksheets.Item($idx).Range("A1").Select() | Out-Null
}

# --- icons sheet: add a new header row (A1/B1) above the existing data, matching the
#     style used for the same header elsewhere in the workbook, then move the selection
#     and make sure this sheet stays the active one. ---
$wsIcons = $wb.Worksheets.Item(5)
$wsFieldnames = $wb.Worksheets.Item(1)

$wsFieldnames.Range("A1:B1").Copy() | Out-Null
$wsIcons.Range("A1:B1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$wsIcons.Range("A1").Value = "my OS-drive"
$wsIcons.Range("B1").Value = "my data-drive"

$wsIcons.Activate() | Out-Null
$wsIcons.Range("B21").Select() | Out-Null
